$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values that look numeric must be forced to Text format so Excel
# keeps the exact original string (e.g. trailing zeros) instead of coercing
# them into a floating point number.
$textCells = @('D5', 'D6', 'D8', 'D10', 'D11', 'D14', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D29', 'D30', 'D32', 'D33', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D43', 'D44', 'D46', 'D47', 'D49', 'D51')
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '61.485.75'
$ws.Range('E2').Value = '  -2.80%  '
$ws.Range('D3').Value = '2.891.73'
$ws.Range('E3').Value = '  -2.64%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '570.20'
$ws.Range('E5').Value = '  -4.47%  '
$ws.Range('D6').Value = '142.90'
$ws.Range('E6').Value = '  -4.72%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.507'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = '2.890.00'
$ws.Range('E9').Value = '  -2.66%  '
$ws.Range('D10').Value = '6.74'
$ws.Range('E10').Value = '  -8.34%  '
$ws.Range('D11').Value = '0.146'
$ws.Range('E11').Value = '  -6.00%  '
$ws.Range('E12').Value = '  -3.23%  '
$ws.Range('E13').Value = '  -4.68%  '
$ws.Range('D14').Value = '31.93'
$ws.Range('E14').Value = '  -4.38%  '
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').Value = '3.373.45'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').Value = '61.541.08'
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '6.63'
$ws.Range('E18').Value = '  -2.70%  '
$ws.Range('D19').Value = '2.874.97'
$ws.Range('E19').Value = '  -3.21%  '
$ws.Range('D20').Value = '432.64'
$ws.Range('E20').Value = '  -3.17%  '
$ws.Range('D21').Value = '13.17'
$ws.Range('E21').Value = '  -3.40%  '
$ws.Range('D22').Value = '0.656'
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('D23').Value = '6.88'
$ws.Range('E23').Value = '  -4.08%  '
$ws.Range('D24').Value = '79.45'
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('E25').Value = '  -1.42%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '10.01'
$ws.Range('E27').Value = '  -12.41%  '
$ws.Range('D28').Value = '2.02'
$ws.Range('E28').Value = '  -8.61%  '
$ws.Range('D29').Value = '0.0000104'
$ws.Range('E29').Value = '  -4.30%  '
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  -4.56%  '
$ws.Range('E31').Value = '  -5.57%  '
$ws.Range('D32').Value = '2.06'
$ws.Range('E32').Value = '  -8.49%  '
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('E34').Value = '  -3.99%  '
$ws.Range('D35').Value = '25.62'
$ws.Range('E35').Value = '  -4.75%  '
$ws.Range('D36').Value = '0.956'
$ws.Range('E36').Value = '  -4.32%  '
$ws.Range('D37').Value = '5.42'
$ws.Range('E37').Value = '  -4.95%  '
$ws.Range('D38').Value = '48.79'
$ws.Range('E38').Value = '  -2.19%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').Value = '2.84'
$ws.Range('E39').Value = '  -14.70%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  -7.09%  '
$ws.Range('E41').Value = '  -4.32%  '
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('D43').Value = '39.60'
$ws.Range('E43').Value = '  -4.31%  '
$ws.Range('D44').Value = '0.267'
$ws.Range('E44').Value = '  -6.99%  '
$ws.Range('D45').Value = '2.684.41'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('D46').Value = '133.92'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').Value = '0.0334'
$ws.Range('E47').Value = '  -2.97%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').Value = '339.78'
$ws.Range('E49').Value = '  -8.87%  '
$ws.Range('E50').Value = '  -2.62%  '
$ws.Range('D51').Value = '21.54'
$ws.Range('E51').Value = '  -7.97%  '
